$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '42.903.10'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.63%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.535.89'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.66%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '311.28'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '100.83'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.35%  '

$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -0.85%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '35.82'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.87%  '

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("E12").Value = '  -0.51%  '

$ws.Range("E13").Value = '  +1.80%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.924.99'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.86%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.49'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.21%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.522.15'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.40%  '

$ws.Range("E17").Value = '  -1.87%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '42.870.40'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.36%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.70'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("E20").Value = '  +0.78%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0953'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.11%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '69.82'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.57%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '243.70'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.35%  '

$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("E27").Value = '  -3.94%  '

$ws.Range("E28").Value = '  -2.34%  '

$ws.Range("E29").Value = '  +1.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '38.86'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.59%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '159.79'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.25%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.85'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.43%  '

$ws.Range("E33").Value = '  +7.43%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.61%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0793'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '18.34'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("E37").Value = '  -4.53%  '

$ws.Range("E38").Value = '  -4.62%  '

$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("E40").Value = '  +0.39%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.20'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.64%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '21.89'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.95%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.34'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +5.04%  '

$ws.Range("E44").Value = '  +0.16%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0300'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.006.33'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.05%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.19'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.80%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.777.35'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.92%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.193'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.32%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '80.04'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.01%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '72.41'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '
